$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$rows = @(
    @(375, 44449, 0, 1, 31.25976867771178),
    @(376, 44450, 0, 1, 31.25976867771178),
    @(377, 44451, 1, 2, 62.51953735542357),
    @(378, 44452, 0, 2, 62.51953735542357),
    @(379, 44453, 0, 2, 62.51953735542357),
    @(380, 44454, 0, 2, 62.51953735542357),
    @(381, 44455, 1, 2, 62.51953735542357),
    @(382, 44456, 1, 3, 93.77930603313536),
    @(383, 44457, 0, 3, 93.77930603313536),
    @(384, 44458, 2, 4, 125.0390747108471),
    @(385, 44459, 0, 4, 125.0390747108471)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $dateVal = $r[1]
    $newPos = $r[2]
    $sumMobile = $r[3]
    $sumMobile100k = $r[4]

    # Copy the style (incl. date number format) from the last existing row (374) column A
    $ws.Range("A374").Copy($ws.Range("A$rowNum"))
    $ws.Range("A$rowNum").Value = $dateVal

    $ws.Range("B$rowNum").Value = $newPos
    $ws.Range("C$rowNum").Value = $sumMobile
    $ws.Range("D$rowNum").Value = $sumMobile100k
}
